# Result analysis per class and exported into Excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): A1 already holds 0 with the bold/bordered style.
# Extend the same values/style across B1:G1 ---
$headerValues = @(1,2,3,4,5,6)
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $col = $i + 2  # B..G
    $ws.Cells.Item(1, $col).Value = $headerValues[$i]
}

# Copy A1's formatting (bold font, borders, centered) onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-8: student id, three module marks, total, average, class ---
$data = @(
    @(1, 75, 88, 74, 238, 79, "First class"),
    @(2, 85, 90, 69, 246, 82, "First class With Distinction"),
    @(3, 65, 78, 99, 245, 81, "First class With Distinction"),
    @(4, 92, 74, 45, 215, 71, "First class"),
    @(5, 75, 96, 74, 250, 83, "First class With Distinction"),
    @(6, 55, 47, 36, 144, 48, "Fail"),
    @(7, 68, 45, 87, 207, 69, "First class")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}

Write-Output "Populated results sheet A1:G8"
